$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Disease Ontology (DO) source_version: v2023-08-08 -> v2023-09-29
$ws.Range("E3").Value = "v2023-09-29"

# Experimental Factor Ontology (EFO) source_version: v3.57.0 -> v3.58.0
$ws.Range("E4").Value = "v3.58.0"

# Move the active cell selection to E5, matching the saved cursor position
$ws.Range("E5").Select()
